# Auto-generated: update computed market-price columns (H-N) across all sheets
# per the scheduled runner's refreshed data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(20, 8).Value = 4666.6665  # H20: 4175 -> 4666.6665
$ws.Cells.Item(20, 9).Value = 4000  # I20: 3900 -> 4000
$ws.Cells.Item(20, 11).Value = 4000  # K20: 3900 -> 4000
$ws.Cells.Item(20, 13).Value = -3770  # M20: -3670 -> -3770
$ws.Cells.Item(35, 8).Value = 4666.6665  # H35: 4175 -> 4666.6665
$ws.Cells.Item(35, 9).Value = 4000  # I35: 3900 -> 4000
$ws.Cells.Item(35, 11).Value = 4000  # K35: 3900 -> 4000
$ws.Cells.Item(35, 13).Value = -3621  # M35: -3521 -> -3621
$ws.Cells.Item(38, 8).Value = 64111.125  # H38: 2014.1666 -> 64111.125
$ws.Cells.Item(38, 9).Value = 72555.86  # I38: 2014.1666 -> 72555.86
$ws.Cells.Item(38, 10).Value = 4998  # J38: 0 -> 4998
$ws.Cells.Item(38, 11).Value = 217667.58  # K38: 6042.4998 -> 217667.58
$ws.Cells.Item(38, 12).Value = 14994  # L38: 0 -> 14994
$ws.Cells.Item(38, 13).Value = -217295.58  # M38: -5670.4998 -> -217295.58
$ws.Cells.Item(38, 14).Value = -15738  # N38: None -> -15738
$ws.Cells.Item(98, 8).Value = 1789.4138  # H98: 1738.6666 -> 1789.4138
$ws.Cells.Item(98, 9).Value = 1377.28  # I98: 1334.5769 -> 1377.28
$ws.Cells.Item(98, 11).Value = 1377.28  # K98: 1334.5769 -> 1377.28
$ws.Cells.Item(98, 13).Value = 120.72  # M98: 163.4231 -> 120.72
$ws.Cells.Item(112, 8).Value = 3005.3  # H112: 3219.3333 -> 3005.3
$ws.Cells.Item(112, 10).Value = 3232.8235  # J112: 3469.889 -> 3232.8235
$ws.Cells.Item(112, 12).Value = 9698.470499999999  # L112: 10409.667 -> 9698.470499999999
$ws.Cells.Item(112, 14).Value = -11914.4705  # N112: -12625.667 -> -11914.4705
$ws.Cells.Item(122, 8).Value = 1789.4138  # H122: 1738.6666 -> 1789.4138
$ws.Cells.Item(122, 9).Value = 1377.28  # I122: 1334.5769 -> 1377.28
$ws.Cells.Item(122, 11).Value = 4131.84  # K122: 4003.7307 -> 4131.84
$ws.Cells.Item(122, 13).Value = -1681.84  # M122: -1553.7307 -> -1681.84
$ws.Cells.Item(138, 8).Value = 5160.6987  # H138: 5224.6978 -> 5160.6987
$ws.Cells.Item(138, 9).Value = 2482.4062  # I138: 2790.4644 -> 2482.4062
$ws.Cells.Item(138, 10).Value = 6841.1963  # J138: 6399.8447 -> 6841.1963
$ws.Cells.Item(138, 11).Value = 7447.2186  # K138: 8371.393199999999 -> 7447.2186
$ws.Cells.Item(138, 12).Value = 20523.5889  # L138: 19199.5341 -> 20523.5889
$ws.Cells.Item(138, 13).Value = -2307.2186  # M138: -3231.393199999999 -> -2307.2186
$ws.Cells.Item(138, 14).Value = -30803.5889  # N138: -29479.5341 -> -30803.5889

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 108.5  # H4: 121.28571 -> 108.5
$ws.Cells.Item(4, 9).Value = 90.2  # I4: 108.166664 -> 90.2
$ws.Cells.Item(4, 11).Value = 90.2  # K4: 108.166664 -> 90.2
$ws.Cells.Item(4, 13).Value = 25.8  # M4: 7.833336000000003 -> 25.8
$ws.Cells.Item(32, 8).Value = 2585.9775  # H32: 2611.5557 -> 2585.9775
$ws.Cells.Item(32, 9).Value = 1761.7882  # I32: 1798.1395 -> 1761.7882
$ws.Cells.Item(32, 11).Value = 1761.7882  # K32: 1798.1395 -> 1761.7882
$ws.Cells.Item(32, 13).Value = -1474.7882  # M32: -1511.1395 -> -1474.7882
$ws.Cells.Item(61, 8).Value = 2981.257  # H61: 3754.7273 -> 2981.257
$ws.Cells.Item(61, 9).Value = 2735.9355  # I61: 3634.45 -> 2735.9355
$ws.Cells.Item(61, 10).Value = 4882.5  # J61: 4957.5 -> 4882.5
$ws.Cells.Item(61, 11).Value = 2735.9355  # K61: 3634.45 -> 2735.9355
$ws.Cells.Item(61, 12).Value = 4882.5  # L61: 4957.5 -> 4882.5
$ws.Cells.Item(61, 13).Value = -2523.9355  # M61: -3422.45 -> -2523.9355
$ws.Cells.Item(61, 14).Value = -5306.5  # N61: -5381.5 -> -5306.5
$ws.Cells.Item(74, 8).Value = 1532.2094  # H74: 1699.5278 -> 1532.2094
$ws.Cells.Item(74, 9).Value = 1281.3784  # I74: 1380.75 -> 1281.3784
$ws.Cells.Item(74, 10).Value = 3079  # J74: 4249.75 -> 3079
$ws.Cells.Item(74, 11).Value = 1281.3784  # K74: 1380.75 -> 1281.3784
$ws.Cells.Item(74, 12).Value = 3079  # L74: 4249.75 -> 3079
$ws.Cells.Item(74, 13).Value = -407.3784000000001  # M74: -506.75 -> -407.3784000000001
$ws.Cells.Item(74, 14).Value = -4827  # N74: -5997.75 -> -4827
$ws.Cells.Item(77, 8).Value = 1532.2094  # H77: 1699.5278 -> 1532.2094
$ws.Cells.Item(77, 9).Value = 1281.3784  # I77: 1380.75 -> 1281.3784
$ws.Cells.Item(77, 10).Value = 3079  # J77: 4249.75 -> 3079
$ws.Cells.Item(77, 11).Value = 6406.892  # K77: 6903.75 -> 6406.892
$ws.Cells.Item(77, 12).Value = 15395  # L77: 21248.75 -> 15395
$ws.Cells.Item(77, 13).Value = -2038.892  # M77: -2535.75 -> -2038.892
$ws.Cells.Item(77, 14).Value = -24131  # N77: -29984.75 -> -24131
$ws.Cells.Item(132, 8).Value = 1112.0667  # H132: 1447.6364 -> 1112.0667
$ws.Cells.Item(132, 9).Value = 250.56522  # I132: 283.26666 -> 250.56522
$ws.Cells.Item(132, 11).Value = 751.6956600000001  # K132: 849.79998 -> 751.6956600000001
$ws.Cells.Item(132, 13).Value = 1778.30434  # M132: 1680.20002 -> 1778.30434
$ws.Cells.Item(136, 8).Value = 2981.257  # H136: 3754.7273 -> 2981.257
$ws.Cells.Item(136, 9).Value = 2735.9355  # I136: 3634.45 -> 2735.9355
$ws.Cells.Item(136, 10).Value = 4882.5  # J136: 4957.5 -> 4882.5
$ws.Cells.Item(136, 11).Value = 8207.806500000001  # K136: 10903.35 -> 8207.806500000001
$ws.Cells.Item(136, 12).Value = 14647.5  # L136: 14872.5 -> 14647.5
$ws.Cells.Item(136, 13).Value = -5657.806500000001  # M136: -8353.349999999999 -> -5657.806500000001
$ws.Cells.Item(136, 14).Value = -19747.5  # N136: -19972.5 -> -19747.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 200004000  # H20: 38463588 -> 200004000
$ws.Cells.Item(20, 9).Value = 1000000000  # I20: 71429780 -> 1000000000
$ws.Cells.Item(20, 10).Value = 5003.25  # J20: 3040.9167 -> 5003.25
$ws.Cells.Item(20, 11).Value = 1000000000  # K20: 71429780 -> 1000000000
$ws.Cells.Item(20, 12).Value = 5003.25  # L20: 3040.9167 -> 5003.25
$ws.Cells.Item(20, 13).Value = -999999753  # M20: -71429533 -> -999999753
$ws.Cells.Item(20, 14).Value = -5497.25  # N20: -3534.9167 -> -5497.25
$ws.Cells.Item(36, 8).Value = 6000  # H36: 5666.3335 -> 6000
$ws.Cells.Item(36, 9).Value = 6000  # I36: 5666.3335 -> 6000
$ws.Cells.Item(36, 11).Value = 6000  # K36: 5666.3335 -> 6000
$ws.Cells.Item(36, 13).Value = -5466  # M36: -5132.3335 -> -5466
$ws.Cells.Item(37, 8).Value = 2450.4  # H37: 2311.75 -> 2450.4
$ws.Cells.Item(37, 9).Value = 2450.4  # I37: 1645.2 -> 2450.4
$ws.Cells.Item(37, 10).Value = 0  # J37: 2533.9333 -> 0
$ws.Cells.Item(37, 11).Value = 2450.4  # K37: 1645.2 -> 2450.4
$ws.Cells.Item(37, 12).Value = 0  # L37: 2533.9333 -> 0
$ws.Cells.Item(37, 13).Value = -2313.4  # M37: -1508.2 -> -2313.4
$ws.Cells.Item(37, 14).Value = ""  # N37: -2807.9333 -> None
$ws.Cells.Item(38, 8).Value = 9000  # H38: 0 -> 9000
$ws.Cells.Item(38, 10).Value = 9000  # J38: 0 -> 9000
$ws.Cells.Item(38, 12).Value = 9000  # L38: 0 -> 9000
$ws.Cells.Item(38, 14).Value = -9832  # N38: None -> -9832
$ws.Cells.Item(39, 8).Value = 54324  # H39: 0 -> 54324
$ws.Cells.Item(39, 9).Value = 100048  # I39: 0 -> 100048
$ws.Cells.Item(39, 10).Value = 8600  # J39: 0 -> 8600
$ws.Cells.Item(39, 11).Value = 100048  # K39: 0 -> 100048
$ws.Cells.Item(39, 12).Value = 8600  # L39: 0 -> 8600
$ws.Cells.Item(39, 13).Value = -99659  # M39: None -> -99659
$ws.Cells.Item(39, 14).Value = -9378  # N39: None -> -9378
$ws.Cells.Item(134, 8).Value = 21661.754  # H134: 18867.87 -> 21661.754
$ws.Cells.Item(134, 9).Value = 2717  # I134: 2386.6545 -> 2717
$ws.Cells.Item(134, 10).Value = 170062.33  # J134: 169945.67 -> 170062.33
$ws.Cells.Item(134, 11).Value = 8151  # K134: 7159.9635 -> 8151
$ws.Cells.Item(134, 12).Value = 510186.99  # L134: 509837.01 -> 510186.99
$ws.Cells.Item(134, 13).Value = -5616  # M134: -4624.9635 -> -5616
$ws.Cells.Item(134, 14).Value = -515256.99  # N134: -514907.01 -> -515256.99

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 436.125  # H7: 423.39395 -> 436.125
$ws.Cells.Item(7, 9).Value = 403.16  # I7: 388.26923 -> 403.16
$ws.Cells.Item(7, 11).Value = 403.16  # K7: 388.26923 -> 403.16
$ws.Cells.Item(7, 13).Value = -290.16  # M7: -275.26923 -> -290.16
$ws.Cells.Item(22, 8).Value = 125  # H22: 123.125 -> 125
$ws.Cells.Item(22, 9).Value = 100  # I22: 99 -> 100
$ws.Cells.Item(22, 10).Value = 150  # J22: 163.33333 -> 150
$ws.Cells.Item(22, 11).Value = 100  # K22: 99 -> 100
$ws.Cells.Item(22, 12).Value = 150  # L22: 163.33333 -> 150
$ws.Cells.Item(22, 13).Value = 250  # M22: 251 -> 250
$ws.Cells.Item(22, 14).Value = -850  # N22: -863.3333299999999 -> -850
$ws.Cells.Item(31, 8).Value = 70560.734  # H31: 96215.73 -> 70560.734
$ws.Cells.Item(31, 9).Value = 1039.8  # I31: 1726.6666 -> 1039.8
$ws.Cells.Item(31, 11).Value = 1039.8  # K31: 1726.6666 -> 1039.8
$ws.Cells.Item(31, 13).Value = -744.8  # M31: -1431.6666 -> -744.8
$ws.Cells.Item(34, 8).Value = 70560.734  # H34: 96215.73 -> 70560.734
$ws.Cells.Item(34, 9).Value = 1039.8  # I34: 1726.6666 -> 1039.8
$ws.Cells.Item(34, 11).Value = 1039.8  # K34: 1726.6666 -> 1039.8
$ws.Cells.Item(34, 13).Value = -837.8  # M34: -1524.6666 -> -837.8
$ws.Cells.Item(58, 8).Value = 181797.48  # H58: 195772.98 -> 181797.48
$ws.Cells.Item(58, 9).Value = 224632.44  # I58: 246512.1 -> 224632.44
$ws.Cells.Item(58, 10).Value = 6563.5454  # J58: 6654.4546 -> 6563.5454
$ws.Cells.Item(58, 11).Value = 224632.44  # K58: 246512.1 -> 224632.44
$ws.Cells.Item(58, 12).Value = 6563.5454  # L58: 6654.4546 -> 6563.5454
$ws.Cells.Item(58, 13).Value = -224429.44  # M58: -246309.1 -> -224429.44
$ws.Cells.Item(58, 14).Value = -6969.5454  # N58: -7060.4546 -> -6969.5454
$ws.Cells.Item(94, 8).Value = 2298  # H94: 2982.6 -> 2298
$ws.Cells.Item(94, 9).Value = 0  # I94: 6999 -> 0
$ws.Cells.Item(94, 10).Value = 2298  # J94: 1978.5 -> 2298
$ws.Cells.Item(94, 11).Value = 0  # K94: 6999 -> 0
$ws.Cells.Item(94, 12).Value = 2298  # L94: 1978.5 -> 2298
$ws.Cells.Item(94, 13).Value = ""  # M94: -6548 -> None
$ws.Cells.Item(94, 14).Value = -3200  # N94: -2880.5 -> -3200
$ws.Cells.Item(134, 8).Value = 573141.7  # H134: 835903.75 -> 573141.7
$ws.Cells.Item(134, 9).Value = 346786.94  # I134: 529170.8 -> 346786.94
$ws.Cells.Item(134, 10).Value = 1667189.6  # J134: 2001488.8 -> 1667189.6
$ws.Cells.Item(134, 11).Value = 1040360.82  # K134: 1587512.4 -> 1040360.82
$ws.Cells.Item(134, 12).Value = 5001568.800000001  # L134: 6004466.4 -> 5001568.800000001
$ws.Cells.Item(134, 13).Value = -1037825.82  # M134: -1584977.4 -> -1037825.82
$ws.Cells.Item(134, 14).Value = -5006638.800000001  # N134: -6009536.4 -> -5006638.800000001
$ws.Cells.Item(136, 8).Value = 181797.48  # H136: 195772.98 -> 181797.48
$ws.Cells.Item(136, 9).Value = 224632.44  # I136: 246512.1 -> 224632.44
$ws.Cells.Item(136, 10).Value = 6563.5454  # J136: 6654.4546 -> 6563.5454
$ws.Cells.Item(136, 11).Value = 673897.3200000001  # K136: 739536.3 -> 673897.3200000001
$ws.Cells.Item(136, 12).Value = 19690.6362  # L136: 19963.3638 -> 19690.6362
$ws.Cells.Item(136, 13).Value = -671347.3200000001  # M136: -736986.3 -> -671347.3200000001
$ws.Cells.Item(136, 14).Value = -24790.6362  # N136: -25063.3638 -> -24790.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 908623.7  # H5: 788784.3 -> 908623.7
$ws.Cells.Item(5, 9).Value = 160540.6  # I5: 133922.17 -> 160540.6
$ws.Cells.Item(5, 10).Value = 1376175.6  # J5: 1225359.1 -> 1376175.6
$ws.Cells.Item(5, 11).Value = 481621.8  # K5: 401766.51 -> 481621.8
$ws.Cells.Item(5, 12).Value = 4128526.8  # L5: 3676077.3 -> 4128526.8
$ws.Cells.Item(5, 13).Value = -481509.8  # M5: -401654.51 -> -481509.8
$ws.Cells.Item(5, 14).Value = -4128750.8  # N5: -3676301.3 -> -4128750.8
$ws.Cells.Item(32, 8).Value = 6001050  # H32: 9000700 -> 6001050
$ws.Cells.Item(32, 10).Value = 1201259.8  # J32: 2000933.4 -> 1201259.8
$ws.Cells.Item(32, 12).Value = 3603779.4  # L32: 6002800.199999999 -> 3603779.4
$ws.Cells.Item(32, 14).Value = -3604345.4  # N32: -6003366.199999999 -> -3604345.4
$ws.Cells.Item(113, 8).Value = 1684690.6  # H113: 2059046.4 -> 1684690.6
$ws.Cells.Item(113, 9).Value = 9259696  # I113: 9259784 -> 9259696
$ws.Cells.Item(113, 10).Value = 1356  # J113: 1692.7858 -> 1356
$ws.Cells.Item(113, 11).Value = 27779088  # K113: 27779352 -> 27779088
$ws.Cells.Item(113, 12).Value = 4068  # L113: 5078.357400000001 -> 4068
$ws.Cells.Item(113, 13).Value = -27776918  # M113: -27777182 -> -27776918
$ws.Cells.Item(113, 14).Value = -8408  # N113: -9418.357400000001 -> -8408
$ws.Cells.Item(135, 8).Value = 908623.7  # H135: 788784.3 -> 908623.7
$ws.Cells.Item(135, 9).Value = 160540.6  # I135: 133922.17 -> 160540.6
$ws.Cells.Item(135, 10).Value = 1376175.6  # J135: 1225359.1 -> 1376175.6
$ws.Cells.Item(135, 11).Value = 1444865.4  # K135: 1205299.53 -> 1444865.4
$ws.Cells.Item(135, 12).Value = 12385580.4  # L135: 11028231.9 -> 12385580.4
$ws.Cells.Item(135, 13).Value = -1442330.4  # M135: -1202764.53 -> -1442330.4
$ws.Cells.Item(135, 14).Value = -12390650.4  # N135: -11033301.9 -> -12390650.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1037.7273  # H97: 1077.8572 -> 1037.7273
$ws.Cells.Item(97, 9).Value = 829.9286  # I97: 878.7692 -> 829.9286
$ws.Cells.Item(97, 11).Value = 829.9286  # K97: 878.7692 -> 829.9286
$ws.Cells.Item(97, 13).Value = -333.9286  # M97: -382.7692 -> -333.9286
$ws.Cells.Item(132, 8).Value = 277719.78  # H132: 346859.56 -> 277719.78
$ws.Cells.Item(132, 9).Value = 315502.12  # I132: 360481.34 -> 315502.12
$ws.Cells.Item(132, 10).Value = 126590.375  # J132: 251507 -> 126590.375
$ws.Cells.Item(132, 11).Value = 946506.36  # K132: 1081444.02 -> 946506.36
$ws.Cells.Item(132, 12).Value = 379771.125  # L132: 754521 -> 379771.125
$ws.Cells.Item(132, 13).Value = -943976.36  # M132: -1078914.02 -> -943976.36
$ws.Cells.Item(132, 14).Value = -384831.125  # N132: -759581 -> -384831.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 638.5  # H22: 596.5 -> 638.5
$ws.Cells.Item(22, 9).Value = 575  # I22: 528.8 -> 575
$ws.Cells.Item(22, 10).Value = 670.25  # J22: 664.2 -> 670.25
$ws.Cells.Item(22, 11).Value = 575  # K22: 528.8 -> 575
$ws.Cells.Item(22, 12).Value = 670.25  # L22: 664.2 -> 670.25
$ws.Cells.Item(22, 13).Value = -280  # M22: -233.8 -> -280
$ws.Cells.Item(22, 14).Value = -1260.25  # N22: -1254.2 -> -1260.25
$ws.Cells.Item(27, 8).Value = 638.5  # H27: 596.5 -> 638.5
$ws.Cells.Item(27, 9).Value = 575  # I27: 528.8 -> 575
$ws.Cells.Item(27, 10).Value = 670.25  # J27: 664.2 -> 670.25
$ws.Cells.Item(27, 11).Value = 575  # K27: 528.8 -> 575
$ws.Cells.Item(27, 12).Value = 670.25  # L27: 664.2 -> 670.25
$ws.Cells.Item(27, 13).Value = -468  # M27: -421.8 -> -468
$ws.Cells.Item(27, 14).Value = -884.25  # N27: -878.2 -> -884.25
$ws.Cells.Item(40, 8).Value = 146009.72  # H40: 57365.055 -> 146009.72
$ws.Cells.Item(40, 9).Value = 169732.17  # I40: 68517.734 -> 169732.17
$ws.Cells.Item(40, 10).Value = 3675  # J40: 1601.6666 -> 3675
$ws.Cells.Item(40, 11).Value = 169732.17  # K40: 68517.734 -> 169732.17
$ws.Cells.Item(40, 12).Value = 3675  # L40: 1601.6666 -> 3675
$ws.Cells.Item(40, 13).Value = -169596.17  # M40: -68381.734 -> -169596.17
$ws.Cells.Item(40, 14).Value = -3947  # N40: -1873.6666 -> -3947
$ws.Cells.Item(46, 8).Value = 3846.8096  # H46: 5079.125 -> 3846.8096
$ws.Cells.Item(46, 9).Value = 3375.7693  # I46: 5368.6875 -> 3375.7693
$ws.Cells.Item(46, 10).Value = 4612.25  # J46: 4500 -> 4612.25
$ws.Cells.Item(46, 11).Value = 3375.7693  # K46: 5368.6875 -> 3375.7693
$ws.Cells.Item(46, 12).Value = 4612.25  # L46: 4500 -> 4612.25
$ws.Cells.Item(46, 13).Value = -3187.7693  # M46: -5180.6875 -> -3187.7693
$ws.Cells.Item(46, 14).Value = -4988.25  # N46: -4876 -> -4988.25
$ws.Cells.Item(55, 8).Value = 1348.8572  # H55: 714.1724 -> 1348.8572
$ws.Cells.Item(55, 9).Value = 228.33333  # I55: 193.5 -> 228.33333
$ws.Cells.Item(55, 10).Value = 3365.8  # J55: 1566.1818 -> 3365.8
$ws.Cells.Item(55, 11).Value = 228.33333  # K55: 193.5 -> 228.33333
$ws.Cells.Item(55, 12).Value = 3365.8  # L55: 1566.1818 -> 3365.8
$ws.Cells.Item(55, 13).Value = -55.33332999999999  # M55: -20.5 -> -55.33332999999999
$ws.Cells.Item(55, 14).Value = -3711.8  # N55: -1912.1818 -> -3711.8
$ws.Cells.Item(59, 8).Value = 0  # H59: 99750 -> 0
$ws.Cells.Item(59, 10).Value = 0  # J59: 99750 -> 0
$ws.Cells.Item(59, 12).Value = 0  # L59: 99750 -> 0
$ws.Cells.Item(59, 14).Value = ""  # N59: -101058 -> None
$ws.Cells.Item(68, 8).Value = 4190  # H68: 5537.5 -> 4190
$ws.Cells.Item(68, 9).Value = 3050  # I68: 4466.6665 -> 3050
$ws.Cells.Item(68, 11).Value = 3050  # K68: 4466.6665 -> 3050
$ws.Cells.Item(68, 13).Value = -2301  # M68: -3717.6665 -> -2301
$ws.Cells.Item(71, 8).Value = 4190  # H71: 5537.5 -> 4190
$ws.Cells.Item(71, 9).Value = 3050  # I71: 4466.6665 -> 3050
$ws.Cells.Item(71, 11).Value = 15250  # K71: 22333.3325 -> 15250
$ws.Cells.Item(71, 13).Value = -11506  # M71: -18589.3325 -> -11506
$ws.Cells.Item(136, 8).Value = 575311.3  # H136: 610191.1 -> 575311.3
$ws.Cells.Item(136, 9).Value = 692513.4399999999  # I136: 772338.4399999999 -> 692513.4399999999
$ws.Cells.Item(136, 10).Value = 8834.666999999999  # J136: 7929.7144 -> 8834.666999999999
$ws.Cells.Item(136, 11).Value = 2077540.32  # K136: 2317015.32 -> 2077540.32
$ws.Cells.Item(136, 12).Value = 26504.001  # L136: 23789.1432 -> 26504.001
$ws.Cells.Item(136, 13).Value = -2074990.32  # M136: -2314465.32 -> -2074990.32
$ws.Cells.Item(136, 14).Value = -31604.001  # N136: -28889.1432 -> -31604.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 670.4828  # H113: 709.8889 -> 670.4828
$ws.Cells.Item(113, 9).Value = 616.619  # I113: 666.9474 -> 616.619
$ws.Cells.Item(113, 11).Value = 1849.857  # K113: 2000.8422 -> 1849.857
$ws.Cells.Item(113, 13).Value = 320.143  # M113: 169.1578 -> 320.143
$ws.Cells.Item(136, 8).Value = 113743.61  # H136: 146217.78 -> 113743.61
$ws.Cells.Item(136, 9).Value = 84.8  # I136: 88 -> 84.8
$ws.Cells.Item(136, 11).Value = 254.4  # K136: 264 -> 254.4
$ws.Cells.Item(136, 13).Value = 2295.6  # M136: 2286 -> 2295.6
